$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("D1").Value = "TB người"
$ws.Range("E1").Value = "Tổng frame"
$ws.Range("F1").Value = "Tổng người"

# Update row 2 values
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "1.0m51.0s"
$ws.Range("E2").Value = 370
$ws.Range("F2").Value = 656

# Copy style from D1 to E1/F1 so the new headers match the existing header formatting
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Remove rows 3 and 4 entirely
$ws.Rows.Item(3).Resize(2).Delete()
